$d = $word.ActiveDocument

$replacements = @(
    @("380×2=", "909×7="),
    @("655×2=", "887×4="),
    @("120×6=", "128×7="),
    @("398×2=", "617×8="),
    @("399×6=", "559×4="),
    @("246×6=", "150×2="),
    @("178×8=", "584×7="),
    @("144×5=", "995×3="),
    @("910×5=", "646×8="),
    @("538×7=", "490×3="),
    @("212×7=", "938×2="),
    @("528×2=", "591×5="),
    @("966×8=", "230×6="),
    @("657×5=", "540×2="),
    @("587×4=", "598×6="),
    @("190×6=", "805×5="),
    @("520×8=", "226×6="),
    @("759×3=", "445×5="),
    @("172×4=", "478×6="),
    @("575×6=", "718×8="),
    @("926×6=", "499×8="),
    @("782×2=", "609×5="),
    @("684×4=", "228×9="),
    @("277×7=", "580×5="),
    @("448×8=", "460×8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
